$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns X (24) and AD (30) from 19 to 20 (stored OOXML width units).
# ColumnWidth (characters) = stored width - 0.8333333333333334
$ws.Columns.Item(24).ColumnWidth = 19.166666666666668
$ws.Columns.Item(30).ColumnWidth = 19.166666666666668

# New "Finish Multiplier" marker cell.
$ws.Range("AE2").Value = "UNKNOWN"

# Apply the finish multiplier to the "Total Cost ($)" row values.
$ws.Range("X4").Value = 27933389.2
$ws.Range("Y4").Value = 20950041.9
$ws.Range("AA4").Value = 49700099.4
$ws.Range("AB4").Value = 26350052.7
$ws.Range("AC4").Value = 22200044.4
$ws.Range("AD4").Value = 13200026.4
